$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8898788690567017
$ws.Range("B1").Value = 1.730800747871399
$ws.Range("C1").Value = 4.12269926071167
$ws.Range("D1").Value = 3.60997200012207
$ws.Range("E1").Value = 1.041465759277344
